# Delete unused CBM thickness data rows (lung, nerve, ear, etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBM thickness")

# Rows (1-based, as they currently stand in the sheet) that correspond to the
# now-unused CBM thickness entries being removed:
#   19 - Carlson et al., 2003 (11 mo. FVB mice & Pulmonary alveolus)
#   21 - Carlson et al., 2003 (11 mo. FVB mice & Pancreas)
#   22 - Carlson et al., 2003 (11 mo. FVB mice & Choroid)
#   25 - Carlson et al., 2003 (11 mo. FVB mice & Peripheral nerve)
#   44 - Fraselle-Jacobs et al., 1987 (6 mo. Wistar rat & Adipose)
#   77 - Smith et al., 1995 (6 mo. Sprague-Dawley rat & Inner ear)
#   86 - Chakrabarti et al., 1991 (6 mo. BB rat & Endoneurial capillary)
#
# Delete from the bottom up so earlier row numbers stay valid as we go.
$rowsToDelete = @(86, 77, 44, 25, 22, 21, 19)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

$ws.Activate()
